# Phân công lại người phụ trách (tái gán nhiệm vụ), bỏ "Ngọc" khỏi danh sách
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G8").Value  = "Lan"
$ws.Range("G9").Value  = "Lan"
$ws.Range("G10").Value = "Lan"
$ws.Range("G11").Value = "Huyền"
$ws.Range("G13").Value = "Huyền"
$ws.Range("G14").Value = "Huyền"
$ws.Range("G15").Value = "Độ"
$ws.Range("G16").Value = "Độ"
$ws.Range("G17").Value = "Độ"

$ws.Activate()
$ws.Range("E17").Select()
